# Add a new "break_on_off" column (L) to the 2D training schedule sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column
$ws.Range("L1").Value = "break_on_off"

# break_on_off is 1 only for the trials that are followed by a break screen,
# 0 otherwise. Rows 2..73 correspond to trials 1..72.
$breakRows = @(19, 37, 54)

for ($r = 2; $r -le 73; $r++) {
    if ($breakRows -contains $r) {
        $ws.Cells.Item($r, 12).Value = 1
    } else {
        $ws.Cells.Item($r, 12).Value = 0
    }
}

# Match the author's final selection/highlight of the newly added column.
$ws.Range("L1:L73").Select()
